$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether it is a numeric-looking
# "Price" column value that must be forced to Text so Excel does not
# auto-convert it into a floating point number (losing the exact display).
$updates = @(
    @{ Cell = "D2"; Value = "62.909.63"; ForceText = 1 },
    @{ Cell = "E2"; Value = "  -5.34%  "; ForceText = 0 },
    @{ Cell = "D3"; Value = "3.114.53"; ForceText = 1 },
    @{ Cell = "E3"; Value = "  -5.94%  "; ForceText = 0 },
    @{ Cell = "E4"; Value = "  +0.00%  "; ForceText = 0 },
    @{ Cell = "D5"; Value = "558.60"; ForceText = 1 },
    @{ Cell = "E5"; Value = "  -4.85%  "; ForceText = 0 },
    @{ Cell = "D6"; Value = "163.04"; ForceText = 1 },
    @{ Cell = "E6"; Value = "  -9.61%  "; ForceText = 0 },
    @{ Cell = "E7"; Value = "  +0.01%  "; ForceText = 0 },
    @{ Cell = "D8"; Value = "0.588"; ForceText = 1 },
    @{ Cell = "E8"; Value = "  -10.09%  "; ForceText = 0 },
    @{ Cell = "D9"; Value = "3.103.22"; ForceText = 1 },
    @{ Cell = "E9"; Value = "  -6.19%  "; ForceText = 0 },
    @{ Cell = "D10"; Value = "6.72"; ForceText = 1 },
    @{ Cell = "E10"; Value = "  -1.69%  "; ForceText = 0 },
    @{ Cell = "E11"; Value = "  -8.26%  "; ForceText = 0 },
    @{ Cell = "D12"; Value = "0.379"; ForceText = 1 },
    @{ Cell = "E12"; Value = "  -5.90%  "; ForceText = 0 },
    @{ Cell = "D13"; Value = "3.651.03"; ForceText = 1 },
    @{ Cell = "E13"; Value = "  -6.02%  "; ForceText = 0 },
    @{ Cell = "E14"; Value = "  -1.77%  "; ForceText = 0 },
    @{ Cell = "D15"; Value = "63.024.48"; ForceText = 1 },
    @{ Cell = "E15"; Value = "  -5.14%  "; ForceText = 0 },
    @{ Cell = "D16"; Value = "24.62"; ForceText = 1 },
    @{ Cell = "E16"; Value = "  -7.56%  "; ForceText = 0 },
    @{ Cell = "D17"; Value = "3.106.21"; ForceText = 1 },
    @{ Cell = "E17"; Value = "  -7.04%  "; ForceText = 0 },
    @{ Cell = "E18"; Value = "  -6.07%  "; ForceText = 0 },
    @{ Cell = "D19"; Value = "405.48"; ForceText = 1 },
    @{ Cell = "E19"; Value = "  -4.50%  "; ForceText = 0 },
    @{ Cell = "D20"; Value = "12.51"; ForceText = 1 },
    @{ Cell = "E20"; Value = "  -4.53%  "; ForceText = 0 },
    @{ Cell = "D21"; Value = "5.17"; ForceText = 1 },
    @{ Cell = "E21"; Value = "  -5.51%  "; ForceText = 0 },
    @{ Cell = "E22"; Value = "  -3.84%  "; ForceText = 0 },
    @{ Cell = "D23"; Value = "0.997"; ForceText = 1 },
    @{ Cell = "E23"; Value = "  -0.29%  "; ForceText = 0 },
    @{ Cell = "D24"; Value = "5.68"; ForceText = 1 },
    @{ Cell = "E24"; Value = "  +0.13%  "; ForceText = 0 },
    @{ Cell = "D25"; Value = "68.61"; ForceText = 1 },
    @{ Cell = "E26"; Value = "  -3.21%  "; ForceText = 0 },
    @{ Cell = "E27"; Value = "  -5.06%  "; ForceText = 0 },
    @{ Cell = "E28"; Value = "  -11.40%  "; ForceText = 0 },
    @{ Cell = "D29"; Value = "8.69"; ForceText = 1 },
    @{ Cell = "E29"; Value = "  -4.41%  "; ForceText = 0 },
    @{ Cell = "E30"; Value = "  +0.00%  "; ForceText = 0 },
    @{ Cell = "D32"; Value = "21.28"; ForceText = 1 },
    @{ Cell = "E32"; Value = "  -4.95%  "; ForceText = 0 },
    @{ Cell = "E33"; Value = "  -7.08%  "; ForceText = 0 },
    @{ Cell = "D34"; Value = "4.85"; ForceText = 1 },
    @{ Cell = "E34"; Value = "  -5.64%  "; ForceText = 0 },
    @{ Cell = "E35"; Value = "  -5.85%  "; ForceText = 0 },
    @{ Cell = "D36"; Value = "154.64"; ForceText = 1 },
    @{ Cell = "E36"; Value = "  -3.49%  "; ForceText = 0 },
    @{ Cell = "E37"; Value = "  -7.35%  "; ForceText = 0 },
    @{ Cell = "D38"; Value = "1.33"; ForceText = 1 },
    @{ Cell = "E38"; Value = "  -6.86%  "; ForceText = 0 },
    @{ Cell = "D39"; Value = "2.725.30"; ForceText = 1 },
    @{ Cell = "E39"; Value = "  -4.57%  "; ForceText = 0 },
    @{ Cell = "D40"; Value = "1.65"; ForceText = 1 },
    @{ Cell = "E40"; Value = "  -8.37%  "; ForceText = 0 },
    @{ Cell = "D41"; Value = "23.66"; ForceText = 1 },
    @{ Cell = "E41"; Value = "  -10.11%  "; ForceText = 0 },
    @{ Cell = "E42"; Value = "  -5.06%  "; ForceText = 0 },
    @{ Cell = "D43"; Value = "38.15"; ForceText = 1 },
    @{ Cell = "E43"; Value = "  -3.88%  "; ForceText = 0 },
    @{ Cell = "D44"; Value = "0.694"; ForceText = 1 },
    @{ Cell = "E44"; Value = "  -7.78%  "; ForceText = 0 },
    @{ Cell = "D45"; Value = "0.0604"; ForceText = 1 },
    @{ Cell = "E45"; Value = "  -8.36%  "; ForceText = 0 },
    @{ Cell = "B46"; Value = "RenderToken"; ForceText = 0 },
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; ForceText = 0 },
    @{ Cell = "D46"; Value = "5.26"; ForceText = 1 },
    @{ Cell = "E46"; Value = "  -10.76%  "; ForceText = 0 },
    @{ Cell = "B47"; Value = "VeChain"; ForceText = 0 },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; ForceText = 0 },
    @{ Cell = "D47"; Value = "0.0255"; ForceText = 1 },
    @{ Cell = "E47"; Value = "  -6.17%  "; ForceText = 0 },
    @{ Cell = "D48"; Value = "283.83"; ForceText = 1 },
    @{ Cell = "E48"; Value = "  -9.04%  "; ForceText = 0 },
    @{ Cell = "D49"; Value = "20.84"; ForceText = 1 },
    @{ Cell = "E49"; Value = "  -9.71%  "; ForceText = 0 },
    @{ Cell = "D50"; Value = "0.999"; ForceText = 1 },
    @{ Cell = "E50"; Value = "  +0.00%  "; ForceText = 0 },
    @{ Cell = "D51"; Value = "0.0970"; ForceText = 1 },
    @{ Cell = "E51"; Value = "  -6.72%  "; ForceText = 0 }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText -eq 1) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
